$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Name" row (row 4): the Value column (B4) was empty; fill in the
# generated ValueSet name.
$ws.Range("B4").Value = "TypeactiviteliberaleVs"

# "Date" row (row 8): bump the generation timestamp to match the new
# fsh-generation run.
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
